$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -8.065
$ws.Range("B7").Value = 5.762
$ws.Range("E7").Value = 16.268
$ws.Range("A9").Value = -21.396
$ws.Range("E10").Value = 16.492
$ws.Range("B12").Value = 5.927
$ws.Range("A13").Value = -22.116
$ws.Range("E13").Value = 16.485
$ws.Range("B14").Value = 6.040000000000001
$ws.Range("D15").Value = -8.202000000000002
$ws.Range("A16").Value = -21.988
$ws.Range("E16").Value = 16.737
$ws.Range("A18").Value = -22.055
$ws.Range("B19").Value = 8.425000000000001
$ws.Range("A20").Value = -20.32
$ws.Range("E20").Value = 16.598
$ws.Range("E24").Value = 17.044
$ws.Range("A26").Value = -21.386
$ws.Range("B26").Value = 6.244
$ws.Range("A27").Value = -21.449
$ws.Range("B27").Value = 6.397
$ws.Range("D28").Value = -8.205
$ws.Range("A29").Value = -21.563
$ws.Range("B29").Value = 6.509
$ws.Range("E32").Value = 16.457
$ws.Range("D33").Value = -7.773000000000001
$ws.Range("A35").Value = -19.873
$ws.Range("D35").Value = -7.551
$ws.Range("A36").Value = -21.095
$ws.Range("B37").Value = 8.847000000000001
$ws.Range("B38").Value = 6.273000000000001
$ws.Range("D38").Value = -8.782999999999999
$ws.Range("E39").Value = 16.247
$ws.Range("D43").Value = -7.641
$ws.Range("D44").Value = -7.465999999999999
$ws.Range("A45").Value = -21.704
$ws.Range("D45").Value = -7.565
$ws.Range("B47").Value = 6.004
$ws.Range("D47").Value = -8.125999999999999
$ws.Range("E47").Value = 16.709
$ws.Range("E48").Value = 17.275
$ws.Range("B51").Value = 5.326000000000001
$ws.Range("D51").Value = -8.419999999999998
$ws.Range("B52").Value = 5.083
$ws.Range("E52").Value = 16.441
$ws.Range("D54").Value = -8.183
$ws.Range("A55").Value = -21.682
$ws.Range("B55").Value = 6.206
$ws.Range("E56").Value = 16.759
$ws.Range("A57").Value = -22.186
$ws.Range("D57").Value = -8.077
$ws.Range("D62").Value = -8.071999999999999
$ws.Range("D63").Value = -7.337000000000001
$ws.Range("D67").Value = -6.854000000000001
$ws.Range("A69").Value = -21.459
$ws.Range("B69").Value = 6.16
$ws.Range("B70").Value = 5.667000000000001
$ws.Range("D70").Value = -6.976000000000001
$ws.Range("A76").Value = -22.022
$ws.Range("B76").Value = 5.234
$ws.Range("A78").Value = -20.175
$ws.Range("B81").Value = 6.197000000000001
$ws.Range("D81").Value = -7.414
$ws.Range("A82").Value = -22.143
$ws.Range("A83").Value = -20.094
$ws.Range("B83").Value = 7.536
$ws.Range("E84").Value = 16.617
$ws.Range("D88").Value = -7.963000000000001
$ws.Range("A93").Value = -21.592
$ws.Range("B94").Value = 6.970000000000001
$ws.Range("D96").Value = -7.385
$ws.Range("A97").Value = -22.114
$ws.Range("D99").Value = -8.101000000000001
$ws.Range("B100").Value = 4.997
$ws.Range("E100").Value = 16.347
$ws.Range("E101").Value = 16.793
$ws.Range("B102").Value = 7.13
